# Complete rebuild of database (db name: cmms2)
# Applies the changes described by the OOXML diff to PostgreSQL-persons.xlsx:
#  - adds a new "insert_persons" worksheet with INSERT statements
#  - updates the person_category_type external-reference formula on
#    "create-persons" (sheet1), and its selection
#  - makes "create-accounts" (sheet2) the active/selected sheet

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # create-persons
$ws2 = $wb.Worksheets.Item(2)   # create-accounts

# --- 1. create-persons: fix up the external formula that now points at
#        list!A2 instead of list!$A$3 (row shifted up in the source book) ---
$ws1.Range("B9").Formula = "=[3]list!A2"

# --- 2. add the new "insert_persons" sheet after the existing two sheets ---
$newSheet = $wb.Worksheets.Add($null, $ws2)
$newSheet.Name = "insert_persons"

$newSheet.Range("A1").Value = "INSERT INTO persons VALUES (1, 'hzlopes@senado.leg.br', 'Henrique', 'Zaidan', '2339', 'SEPLAG', null, 'E');"
$newSheet.Range("A2").Value = "INSERT INTO persons VALUES (2, 'pedrohs@senado.leg.br', 'Pedro Henrique', 'Serafim', '2339', 'SEPLAG', null, 'E');"
$newSheet.Range("A3").Value = "INSERT INTO persons VALUES (3, 'igorlima@senado.leg.br', 'Igor', 'Grimaldi', '3629', 'SEGEEN', null, 'E');"
$newSheet.Range("A4").Value = "INSERT INTO persons VALUES (999, 'chefedegabinete@senado.leg.br', 'Chefe', 'de Gabinete', '9999', 'SF', null, 'C');"
$newSheet.Range("A5").Value = "INSERT INTO private.accounts VALUES (1, crypt('123456', gen_salt('bf', 10)), '2019-08-01', '2019-08-02');"
$newSheet.Range("A6").Value = "INSERT INTO private.accounts VALUES (2, crypt('123456', gen_salt('bf', 10)), '2019-08-01', '2019-08-02');"
$newSheet.Range("A7").Value = "INSERT INTO private.accounts VALUES (3, crypt('123456', gen_salt('bf', 10)), '2019-08-01', '2019-08-02');"

$newSheet.Columns.Item(1).ColumnWidth = 95.7109375

$newSheet.PageSetup.LeftMargin = 36.850393728
$newSheet.PageSetup.RightMargin = 36.850393728
$newSheet.PageSetup.TopMargin = 56.692913399999995
$newSheet.PageSetup.BottomMargin = 56.692913399999995
$newSheet.PageSetup.HeaderMargin = 22.67716464
$newSheet.PageSetup.FooterMargin = 22.67716464

# view state for the new sheet: zoom 160%, selection resting on A8
$newSheet.Select()
$excel.ActiveWindow.Zoom = 160
$newSheet.Range("A8").Select()

# --- 3. create-persons: move its own selection to B10 (no longer the active tab) ---
$ws1.Select()
$ws1.Range("B10").Select()

# --- 4. create-accounts becomes the active/selected sheet, selection stays A4 ---
$ws2.Select()
$ws2.Range("A4").Select()

Write-Output "edit complete"
